$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 through 18 (data rows) with the new values.
# These columns hold text values (quantity as text, order number with
# leading zeros, and a literal date/time string), so force text format
# to avoid Excel auto-converting them to numbers/dates.
$ws.Range("B2:B18").NumberFormat = "@"
$ws.Range("D2:D18").NumberFormat = "@"
$ws.Range("G2:G18").NumberFormat = "@"

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 2).Value = "3"
    $ws.Cells.Item($r, 4).Value = "000001702"
    $ws.Cells.Item($r, 7).Value = "10:31:47 2024-05-16"
}
